$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update row 3 (Login) Runmode from YES to NO
$ws.Range("C3").Value = "NO"

# Add new row 5: Complete Course scenario
$ws.Range("A5").Value = "Complete Course"
$ws.Range("B5").Value = "Complete course description"
$ws.Range("C5").Value = "NO"

# Update the active selection to match the authored state
$ws.Range("B14").Select()
